$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Kiwi"
$ws.Range("B14").Value = 69
$ws.Range("C14").Value = "All"

[void]$ws.Range("C14").Select()
